# 設計メモ.xlsx edit: rename Sheet1 -> 構想, insert new "必要な技術" sheet
# after it (before 画面), and add new rows/strings to 構想 and the new sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Rename first sheet (Sheet1) to 構想 ---
$wsKoso = $wb.Worksheets.Item(1)
$wsKoso.Name = "構想"

# --- 2) Insert a brand-new worksheet "必要な技術" right after 構想 (before 画面) ---
# Copying the 画面 sheet (instead of Worksheets.Add()) keeps the workbook's
# usual sheet-level formatting defaults (row height, phoneticPr, namespaces,
# …) instead of the bare-bones blank-sheet defaults .Add() produces.
$wsGamen = $wb.Worksheets.Item(2)
$wsGamen.Copy($null, $wsKoso)
$wsTech = $wb.Worksheets.Item(2)
$wsTech.Name = "必要な技術"

# Wipe the copied content/formatting completely, then drop the two
# left-over formatted rows (top/bottom thick border rows) so nothing of
# the 画面 layout survives.
$wsTech.Cells.Clear()
$wsTech.Rows("12:12").Delete()
$wsTech.Rows("4:4").Delete()

# --- 3) Populate 必要な技術 sheet ---
$wsTech.Range("B3").Value = "HTMLでガントチャートを作れること"
$wsTech.Range("B5").Value = "作ったガントチャートがドラッグアンドドロップで動かせる"
$wsTech.Range("B7").Value = "ドラッグアンドドロップを検知して、イベントを発火させるMVC側に通知"
$wsTech.Range("K3").Select()

# --- 4) Add new content to 構想 sheet ---
$wsKoso.Range("K9").Value = "バックグラウンドで画面に表示したものの警告の可能性がある場所をあらかじめ計算しておけば早い？"

$wsKoso.Range("A36").Value = "根幹を用意したあとで品目追加・LTをGUIベースで変更できるようにしたい（ここがオリジナル）"
$wsKoso.Range("A36").Font.Bold = $true

$wsKoso.Range("A37").Value = "作成した手配をGUIベースで移動できる機能もほしい"
$wsKoso.Range("B38").Value = "移動したときは従属の手配はすべて再計算"

# --- 5) Update the selection on 構想 to match the saved view (D29); this
# also re-activates the 構想 tab since it is the last sheet touched, matching
# tabSelected="1" staying on 構想 in the saved workbook.
$wsKoso.Select()
$wsKoso.Range("D29").Select()
